$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: A2 becomes text "Sala de Aula", C2 becomes "Ambiente de estudo"
$ws.Range("A2").Value = "Sala de Aula"
$ws.Range("C2").Value = "Ambiente de estudo"

# Remove rows 3 and 4 (Professor / Diretor) entirely
$ws.Range("A3:C4").EntireRow.Delete()
